$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sat right after
#    "MP73010" in the title paragraph.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Turn the empty paragraph before the final blank paragraph into
#    the "Version management" commentary paragraph, re-inserting the
#    "_GoBack" bookmark in the middle of the new sentence (right
#    after "management ").
# ------------------------------------------------------------------
$targetPara = $d.Paragraphs(6).Range
$startPos = $targetPara.Start

$beforeBookmark = "The term version management "
$afterBookmark = "refer to the means of effectively tracking and controlling changes to a group of entities usually files and the information contained within them. "

$ins = $d.Range($startPos, $startPos)
$ins.InsertAfter($beforeBookmark + $afterBookmark)

$bookmarkPos = $startPos + $beforeBookmark.Length
$d.Range($bookmarkPos, $bookmarkPos).Bookmarks.Add("_GoBack")

# ------------------------------------------------------------------
# 3. Register the new custom character style used while editing.
# ------------------------------------------------------------------
$newStyle = $d.Styles.Add("e24kjd", 2)
$newStyle.BaseStyle = "DefaultParagraphFont"
